# Q3 Update - 2025
# Applies the data refresh to the UN-MLW (Malawi) sheet:
#  - short-url ("cMy253") changed to "6Z1Ses" across every data row
#  - refreshed refugee/asylum counts for rows 256, 257, 259, 261, 262, 264, 265
#  - "stateless" counts updated for rows 256 and 261
#  - country-of-asylum rows 263-265 shifted to reflect the newly added
#    Zimbabwe entry (Tanzania -> Turkiye -> Uganda -> Zimbabwe)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# short-url column (B) is identical on every data row (2-265); update them all
# in one shot so every cell keeps pointing at the same (now-renamed) value.
$ws.Range("B2:B265").Value = "6Z1Ses"

# Row 256 - Burundi -> Malawi (2024)
$ws.Range("N256").Value = "6662"
$ws.Range("O256").Value = "5751"
$ws.Range("P256").Value = "90"

# Row 257 - Dem. Rep. of the Congo -> Malawi (2024)
$ws.Range("N257").Value = "25649"
$ws.Range("O257").Value = "11128"

# Row 259 - Ethiopia -> Malawi (2024)
$ws.Range("N259").Value = "22"
$ws.Range("O259").Value = "163"

# Row 261 - Rwanda -> Malawi (2024)
$ws.Range("N261").Value = "3102"
$ws.Range("O261").Value = "3950"
$ws.Range("P261").Value = "7"

# Row 262 - Somalia -> Malawi (2024)
$ws.Range("N262").Value = "75"
$ws.Range("O262").Value = "98"

# Row 263 - was United Rep. of Tanzania -> Malawi, now Turkiye -> Malawi
$ws.Range("F263").Value = "196"
$ws.Range("G263").Value = "Türkiye"
$ws.Range("H263").Value = "TUR"
$ws.Range("I263").Value = "TUR"

# Row 264 - was Turkiye -> Malawi, now Uganda -> Malawi
$ws.Range("F264").Value = "199"
$ws.Range("G264").Value = "Uganda"
$ws.Range("H264").Value = "UGA"
$ws.Range("I264").Value = "UGA"
$ws.Range("O264").Value = "9"

# Row 265 - was Uganda -> Malawi, now Zimbabwe -> Malawi
$ws.Range("F265").Value = "214"
$ws.Range("G265").Value = "Zimbabwe"
$ws.Range("H265").Value = "ZIM"
$ws.Range("I265").Value = "ZWE"
$ws.Range("N265").Value = "0"
$ws.Range("O265").Value = "5"
